$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1594223333333333
$ws.Range("H2").Value = 0.478267
$ws.Range("I2").Value = 0.01552338951653915
$ws.Range("J2").Value = 0.01552338951653915
$ws.Range("M2").Value = 24.91851366666667
$ws.Range("N2").Value = 74.75554099999999
$ws.Range("O2").Value = 0.2924799159147552
$ws.Range("P2").Value = 0.2924799159147553
$ws.Range("Q2").Value = 3.972567591938555
$ws.Range("R2").Value = 35.75310832744699
$ws.Range("S2").Value = 0.004540279660509363
$ws.Range("T2").Value = 0.004540279660509364
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1594223333333333
$ws.Range("H3").Value = 0.478267
$ws.Range("I3").Value = 0.01552338951653915
$ws.Range("J3").Value = 0.01552338951653915
$ws.Range("O3").Value = 0.4753125595076708
$ws.Range("P3").Value = 0.4753125595076708
$ws.Range("Q3").Value = 6.455866427737445
$ws.Range("R3").Value = 58.102797849637
$ws.Range("S3").Value = 0.007378462003340767
$ws.Range("T3").Value = 0.007378462003340767
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1594223333333333
$ws.Range("H4").Value = 0.478267
$ws.Range("I4").Value = 0.01552338951653915
$ws.Range("J4").Value = 0.01552338951653915
$ws.Range("M4").Value = 19.78346566666667
$ws.Range("N4").Value = 59.350397
$ws.Range("O4").Value = 0.232207524577574
$ws.Range("P4").Value = 0.232207524577574
$ws.Range("Q4").Value = 3.153926257999889
$ws.Range("R4").Value = 28.385336321999
$ws.Range("S4").Value = 0.003604647852689019
$ws.Range("T4").Value = 0.00360464785268902
$ws.Range("I5").Value = 0.1862883666449807
$ws.Range("J5").Value = 0.1862883666449807
$ws.Range("M5").Value = 24.91851366666667
$ws.Range("N5").Value = 74.75554099999999
$ws.Range("O5").Value = 0.2924799159147552
$ws.Range("P5").Value = 0.2924799159147553
$ws.Range("Q5").Value = 47.67277966584233
$ws.Range("R5").Value = 429.055016992581
$ws.Range("S5").Value = 0.05448560581222104
$ws.Range("T5").Value = 0.05448560581222105
$ws.Range("I6").Value = 0.1862883666449807
$ws.Range("J6").Value = 0.1862883666449807
$ws.Range("O6").Value = 0.4753125595076708
$ws.Range("P6").Value = 0.4753125595076708
$ws.Range("R6").Value = 697.262366403951
$ws.Range("S6").Value = 0.08854520035652917
$ws.Range("T6").Value = 0.08854520035652917
$ws.Range("I7").Value = 0.1862883666449807
$ws.Range("J7").Value = 0.1862883666449807
$ws.Range("M7").Value = 19.78346566666667
$ws.Range("N7").Value = 59.350397
$ws.Range("O7").Value = 0.232207524577574
$ws.Range("P7").Value = 0.232207524577574
$ws.Range("Q7").Value = 37.84867798978634
$ws.Range("R7").Value = 340.638101908077
$ws.Range("S7").Value = 0.04325756047623047
$ws.Range("T7").Value = 0.04325756047623048
$ws.Range("G8").Value = 8.197245333333333
$ws.Range("H8").Value = 24.591736
$ws.Range("I8").Value = 0.7981882438384801
$ws.Range("J8").Value = 0.7981882438384801
$ws.Range("M8").Value = 24.91851366666667
$ws.Range("N8").Value = 74.75554099999999
$ws.Range("O8").Value = 0.2924799159147552
$ws.Range("P8").Value = 0.2924799159147553
$ws.Range("Q8").Value = 204.2631698676862
$ws.Range("R8").Value = 1838.368528809176
$ws.Range("S8").Value = 0.2334540304420248
$ws.Range("T8").Value = 0.2334540304420248
$ws.Range("G9").Value = 8.197245333333333
$ws.Range("H9").Value = 24.591736
$ws.Range("I9").Value = 0.7981882438384801
$ws.Range("J9").Value = 0.7981882438384801
$ws.Range("O9").Value = 0.4753125595076708
$ws.Range("P9").Value = 0.4753125595076708
$ws.Range("Q9").Value = 331.9504854865218
$ws.Range("R9").Value = 2987.554369378695
$ws.Range("S9").Value = 0.3793888971478008
$ws.Range("T9").Value = 0.3793888971478008
$ws.Range("G10").Value = 8.197245333333333
$ws.Range("H10").Value = 24.591736
$ws.Range("I10").Value = 0.7981882438384801
$ws.Range("J10").Value = 0.7981882438384801
$ws.Range("M10").Value = 19.78346566666667
$ws.Range("N10").Value = 59.350397
$ws.Range("O10").Value = 0.232207524577574
$ws.Range("P10").Value = 0.232207524577574
$ws.Range("Q10").Value = 162.1699216132436
$ws.Range("R10").Value = 1459.529294519192
$ws.Range("S10").Value = 0.1853453162486545
$ws.Range("T10").Value = 0.1853453162486545
